$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.008
$ws.Range("G2").Value = -8.987941045109423
$ws.Range("H2").Value = -8.987941045109423
$ws.Range("I2").Value = -28.7863713585846
$ws.Range("J2").Value = -28.7863713585846
$ws.Range("K2").Value = -112.113
$ws.Range("L2").Value = -50.07280035730237
$ws.Range("M2").Value = 35.903
$ws.Range("N2").Value = 0.093184354641958
$ws.Range("O2").Value = -0.3202394013183127
$ws.Range("P2").Value = 33.3
$ws.Range("Q2").Value = 0.08642840457836953
$ws.Range("R2").Value = -0.297021754836638
$ws.Range("S2").Value = 2.602999999999999
$ws.Range("T2").Value = 0.07250090521683422
$ws.Range("U2").Value = 52.023
$ws.Range("V2").Value = 0.1350229697111267
$ws.Range("W2").Value = -0.2351430131362645
$ws.Range("X2").Value = 0.06135765322134987
$ws.Range("Y2").Value = -0.2965006663576144
$ws.Range("Z2").Value = 0.005040165693159735
$ws.Range("AA2").Value = -0.08555866963926978
$ws.Range("AB2").Value = 0.06103582469792974
$ws.Range("AC2").Value = -0.1465774973476638
$ws.Range("AD2").Value = 13.439
$ws.Range("AE2").Value = 0.1034273593546094
$ws.Range("AF2").Value = 13.54242735935461
$ws.Range("AG2").Value = -38.48057264064539
$ws.Range("AH2").Value = 0.03395518124997559
$ws.Range("AI2").Value = 0.03842374624018571
$ws.Range("AJ2").Value = -0.1109559591088418
$ws.Range("AK2").Value = -0.1280865101611614
$ws.Range("AL2").Value = 0.87
$ws.Range("AM2").Value = -24.478
$ws.Range("AN2").Value = -0.4949907918968692
$ws.Range("AO2").Value = -74.17471264367816
$ws.Range("AP2").Value = 1.417332325622298
$ws.Range("AQ2").Value = 2.636326497262848
# Row 3
$ws.Range("K3").Value = -0.411
$ws.Range("U3").Value = 0.235
$ws.Range("V3").Value = 0.04409005628517823
$ws.Range("W3").Value = -3.841121495327103
$ws.Range("X3").Value = 0.06089823974395744
$ws.Range("Y3").Value = -3.90201973507106
$ws.Range("AA3").Value = 0.94919168591224
$ws.Range("AB3").Value = 0.06089823974395744
$ws.Range("AC3").Value = 0.8882934461682825
$ws.Range("AG3").Value = -0.235
$ws.Range("AI3").Value = -0.0
$ws.Range("AJ3").Value = -0.04612365063788027
$ws.Range("AK3").Value = 0.4554263565891473
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
# Row 4
$ws.Range("K4").Value = -0.342
$ws.Range("M4").Value = -0.0
$ws.Range("N4").Value = -0.0
$ws.Range("O4").Value = 0.0
$ws.Range("S4").Value = 0.0
$ws.Range("T4").ClearContents()
$ws.Range("U4").Value = 0.233
$ws.Range("V4").Value = 0.01595890410958904
$ws.Range("W4").Value = -0.01976878612716763
$ws.Range("X4").Value = 0.06519662050344978
$ws.Range("Y4").Value = -0.08496540663061741
$ws.Range("AA4").Value = -0.01422971594265156
$ws.Range("AB4").Value = 0.06237952627189663
$ws.Range("AC4").Value = -0.07660924221454819
$ws.Range("AD4").Value = 1.39
$ws.Range("AF4").Value = 1.39
$ws.Range("AG4").Value = 1.157
$ws.Range("AH4").Value = 0.08692933083176985
$ws.Range("AI4").Value = 0.07358390682901005
$ws.Range("AJ4").Value = 0.07342768293456875
$ws.Range("AK4").Value = 0.06201425738328777
$ws.Range("AL4").Value = 0.078
$ws.Range("AM4").Value = 0.078
$ws.Range("AO4").Value = -3.397435897435898
$ws.Range("AQ4").Value = -3.397435897435898
# Row 5
$ws.Range("D5").Value = -0.106
$ws.Range("G5").Value = -0.4018404907975461
$ws.Range("H5").Value = -0.4018404907975461
$ws.Range("I5").Value = -0.3764284218701305
$ws.Range("J5").Value = -0.3764284218701305
$ws.Range("K5").Value = -6.76
$ws.Range("L5").Value = -4.147239263803681
$ws.Range("U5").Value = 0.513
$ws.Range("V5").Value = 0.06248477466504263
$ws.Range("W5").Value = -0.2119122257053292
$ws.Range("X5").Value = 0.08450572396462933
$ws.Range("Y5").Value = -0.2964179496699585
$ws.Range("Z5").Value = 0.04729038889606876
$ws.Range("AA5").Value = -0.0178014464617719
$ws.Range("AB5").Value = 0.06674899855617444
$ws.Range("AC5").Value = -0.08455044501794634
$ws.Range("AD5").Value = 4.29
$ws.Range("AE5").Value = 0.002891638241563389
$ws.Range("AF5").Value = 4.292891638241564
$ws.Range("AG5").Value = 3.779891638241564
$ws.Range("AH5").Value = 0.3433519030998597
$ws.Range("AI5").Value = 0.158450847386927
$ws.Range("AJ5").Value = 0.315256530441498
$ws.Range("AK5").Value = 0.1422086925592759
$ws.Range("AL5").Value = 0.447
$ws.Range("AM5").Value = -3.633
$ws.Range("AN5").Value = -7.308347529812607
$ws.Range("AO5").Value = -1.378076062639821
$ws.Range("AP5").Value = -6.439338395641506
$ws.Range("AQ5").Value = 0.1695568400770713
# Row 6
$ws.Range("B6").Value = "Falcon Oil & Gas Ltd. (TSXV:FO)"
$ws.Range("D6").Value = -0.167
$ws.Range("G6").Value = -507.4999999999999
$ws.Range("H6").Value = -507.4999999999999
$ws.Range("I6").Value = -572.5
$ws.Range("J6").Value = -572.5
$ws.Range("K6").Value = -1.45
$ws.Range("L6").Value = -362.5
$ws.Range("M6").Value = 0.563
$ws.Range("N6").Value = 0.005040286481647269
$ws.Range("O6").Value = -0.3882758620689655
$ws.Range("S6").Value = 0.563
$ws.Range("U6").Value = 11.5
$ws.Range("V6").Value = 0.1029543419874664
$ws.Range("W6").Value = -0.03295454545454545
$ws.Range("X6").Value = 0.06089823974395744
$ws.Range("Y6").Value = -0.0938527851985029
$ws.Range("Z6").Value = 0.0001298701298701299
$ws.Range("AA6").Value = -0.07435064935064935
$ws.Range("AB6").Value = 0.06089823974395744
$ws.Range("AC6").Value = -0.1352488890946068
$ws.Range("AD6").Value = 0.0
$ws.Range("AE6").Value = 0.0
$ws.Range("AF6").Value = 0.0
$ws.Range("AG6").Value = -11.5
$ws.Range("AH6").Value = 0.0
$ws.Range("AI6").Value = 0.0
$ws.Range("AJ6").Value = -0.1147704590818363
$ws.Range("AK6").Value = -0.362776025236593
$ws.Range("AL6").Value = 0.0
$ws.Range("AM6").Value = -0.057
$ws.Range("AN6").Value = -0.0
$ws.Range("AO6").ClearContents()
$ws.Range("AP6").Value = 5.58252427184466
$ws.Range("AQ6").Value = 40.17543859649123
# Row 7
$ws.Range("B7").Value = "San Leon Energy plc (AIM:SLE)"
$ws.Range("D7").Value = 1.135
$ws.Range("G7").Value = -103.2222222222222
$ws.Range("H7").Value = -103.2222222222222
$ws.Range("I7").Value = -231.1111111111111
$ws.Range("J7").Value = -231.1111111111111
$ws.Range("K7").Value = -52.2
$ws.Range("L7").Value = -580.0
$ws.Range("M7").Value = 35.34
$ws.Range("N7").Value = 0.2367046215673141
$ws.Range("O7").Value = -0.6770114942528734
$ws.Range("P7").Value = 33.3
$ws.Range("Q7").Value = 0.2230408573342264
$ws.Range("R7").Value = -0.6379310344827586
$ws.Range("S7").Value = 2.039999999999999
$ws.Range("T7").Value = 0.05772495755517825
$ws.Range("U7").Value = 35.6
$ws.Range("V7").Value = 0.2384460817146684
$ws.Range("W7").Value = -0.2328278322925959
$ws.Range("X7").Value = 0.06172681830916935
$ws.Range("Y7").Value = -0.2945546506017653
$ws.Range("Z7").Value = 0.000418702023726448
$ws.Range("AA7").Value = -0.09676668992789021
$ws.Range("AB7").Value = 0.0611394156728307
$ws.Range("AC7").Value = -0.1579061056007209
$ws.Range("AD7").Value = 2.74
$ws.Range("AE7").Value = 0.0
$ws.Range("AF7").Value = 2.74
$ws.Range("AG7").Value = -32.86
$ws.Range("AH7").Value = 0.01802157327019206
$ws.Range("AI7").Value = 0.01793898127536991
$ws.Range("AJ7").Value = -0.2822054276880797
$ws.Range("AK7").Value = -0.2805190370496841
$ws.Range("AL7").Value = 0.136
$ws.Range("AM7").Value = -21.064
$ws.Range("AN7").Value = -0.1457446808510638
$ws.Range("AO7").Value = -152.9411764705882
$ws.Range("AP7").Value = 1.747872340425532
$ws.Range("AQ7").Value = 0.9874667679453095
# Row 8
$ws.Range("B8").Value = "Aminex PLC (LSE:AEX)"
$ws.Range("D8").Value = 0.122
$ws.Range("G8").Value = -6.660194174757281
$ws.Range("H8").Value = -6.660194174757281
$ws.Range("I8").Value = -26.01941747572815
$ws.Range("J8").Value = -26.01941747572815
$ws.Range("K8").Value = -14.2
$ws.Range("L8").Value = -27.57281553398058
$ws.Range("U8").Value = 0.932
$ws.Range("V8").Value = 0.03148648648648648
$ws.Range("W8").Value = -0.2374581939799331
$ws.Range("X8").Value = 0.06852466876016737
$ws.Range("Y8").Value = -0.3059828627401005
$ws.Range("Z8").Value = 0.008861739654134046
$ws.Range("AA8").Value = -0.2305773036221286
$ws.Range("AB8").Value = 0.07346933778443522
$ws.Range("AC8").Value = -0.3040466414065638
$ws.Range("AD8").Value = 5.0
$ws.Range("AE8").Value = 0.0
$ws.Range("AF8").Value = 5.0
$ws.Range("AG8").Value = 4.068
$ws.Range("AH8").Value = 0.1445086705202312
$ws.Range("AI8").Value = 0.09823182711198429
$ws.Range("AJ8").Value = 0.1208268979446358
$ws.Range("AK8").Value = 0.0814121037463977
$ws.Range("AL8").Value = 0.056
$ws.Range("AM8").Value = 0.055
$ws.Range("AN8").Value = -2.173913043478261
$ws.Range("AO8").Value = -239.2857142857143
$ws.Range("AP8").Value = -1.768695652173913
$ws.Range("AQ8").Value = -243.6363636363637
# Row 9
$ws.Range("B9").Value = "Providence Resources P.l.c. (ISE:PZQA)"
$ws.Range("K9").Value = -34.4
$ws.Range("U9").Value = 2.55
$ws.Range("V9").Value = 0.04264214046822742
$ws.Range("W9").Value = -0.3940435280641466
$ws.Range("X9").Value = 0.0609884881335304
$ws.Range("Y9").Value = -0.455032016197677
$ws.Range("AA9").Value = -0.2847420962277413
$ws.Range("AB9").Value = 0.06093223372302879
$ws.Range("AC9").Value = -0.3456743299507701
$ws.Range("AD9").Value = 0.019
$ws.Range("AE9").Value = 0.100535721113046
$ws.Range("AF9").Value = 0.119535721113046
$ws.Range("AG9").Value = -2.430464278886954
$ws.Range("AH9").Value = 0.001994937371834923
$ws.Range("AI9").Value = 0.002025358547005386
$ws.Range("AJ9").Value = -0.04236506794654951
$ws.Range("AK9").Value = -0.043040273801618
$ws.Range("AL9").Value = 0.153
$ws.Range("AM9").Value = 0.143
$ws.Range("AN9").Value = -0.005583308845136644
$ws.Range("AO9").Value = -159.4771241830065
$ws.Range("AP9").Value = 0.7142122476893781
$ws.Range("AQ9").Value = -170.6293706293706
# Row 10
$ws.Range("B10").Value = "Petrel Resources Plc (AIM:PET)"
$ws.Range("D10").ClearContents()
$ws.Range("G10").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("I10").ClearContents()
$ws.Range("J10").ClearContents()
$ws.Range("K10").Value = -2.35
$ws.Range("L10").ClearContents()
$ws.Range("U10").Value = 0.46
$ws.Range("V10").Value = 0.06814814814814815
$ws.Range("W10").Value = -0.9437751004016064
$ws.Range("X10").Value = 0.06089823974395744
$ws.Range("Y10").Value = -1.004673340145564
$ws.Range("Z10").Value = 0.0
$ws.Range("AA10").Value = -1.027547004809794
$ws.Range("AB10").Value = 0.06089823974395744
$ws.Range("AC10").Value = -1.088445244553752
$ws.Range("AD10").Value = 0.0
$ws.Range("AE10").Value = 0.0
$ws.Range("AF10").Value = 0.0
$ws.Range("AG10").Value = -0.46
$ws.Range("AH10").Value = 0.0
$ws.Range("AI10").Value = 0.0
$ws.Range("AJ10").Value = -0.07313195548489666
$ws.Range("AK10").Value = -1.074766355140187
$ws.Range("AL10").Value = 0.0
$ws.Range("AM10").Value = 0.0
$ws.Range("AN10").ClearContents()
$ws.Range("AO10").ClearContents()
$ws.Range("AP10").ClearContents()
$ws.Range("AQ10").ClearContents()

Write-Host "Applied all changes"